# "added contingency awareness test"
#
# 1) conditions sheet: bump the duration (column F) of the four main rows -
#    two short blocks (120 -> 180) and two long blocks (300 -> 330). The
#    RAND() driven aux1 column recalculates on its own as a side effect.
# 2) README sheet: build a small contingency-awareness worksheet (minutes
#    per trial, reps, days, totals, and a small summary table) below the
#    existing instructions text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # conditions
$ws2 = $wb.Worksheets.Item(2)   # README

# --- conditions: new block durations (seconds) ---------------------------
$ws1.Range("F2").Value = 180
$ws1.Range("F3").Value = 180
$ws1.Range("F4").Value = 330
$ws1.Range("F5").Value = 330
$ws1.Range("F6").Value = 330
$ws1.Range("F7").Value = 330

# --- README: contingency awareness test table -----------------------------
# Header row + row labels first (this is also the order the new shared
# strings were introduced in, min/trial/rep/u1/u2/o1/o2/days/TOTAL/...).
$ws2.Range("B6").Value = "min"
$ws2.Range("A6").Value = "trial"
$ws2.Range("C6").Value = "rep"

$ws2.Range("A7").Value = "u1"
$ws2.Range("A8").Value = "u2"
$ws2.Range("A9").Value = "o1"
$ws2.Range("A10").Value = "o2"

$ws2.Range("E6").Value = "days"
$ws2.Range("F6").Value = "TOTAL"

$ws2.Range("C14").Value = "total/day"
$ws2.Range("H13").Value = "total 2 days"

$ws2.Range("D6").Value = "total/resp/day"

$ws2.Range("A11").Value = "o1"
$ws2.Range("A12").Value = "o2"
$ws2.Range("I12").Value = "u1"
$ws2.Range("J12").Value = "u2"
$ws2.Range("K12").Value = "o1"
$ws2.Range("L12").Value = "o2"

# row 7 - u1
$ws2.Range("B7").Formula = "=conditions!F2/60"
$ws2.Range("C7").Value = 2
$ws2.Range("D7").Formula = "=B7*C7"
$ws2.Range("E7").Value = 2
$ws2.Range("F7").Formula = "=E7*D7"

# row 8 - u2
$ws2.Range("B8").Formula = "=conditions!F3/60"
$ws2.Range("C8").Value = 2
$ws2.Range("D8").Formula = "=B8*C8"
$ws2.Range("E8").Value = 2
$ws2.Range("F8").Formula = "=E8*D8"

# row 9 - o1
$ws2.Range("B9").Formula = "=conditions!F4/60"
$ws2.Range("C9").Value = 2
$ws2.Range("D9").Formula = "=B9*C9"
$ws2.Range("E9").Value = 2
$ws2.Range("F9").Formula = "=E9*D9"

# row 10 - o2
$ws2.Range("B10").Formula = "=conditions!F5/60"
$ws2.Range("C10").Value = 2
$ws2.Range("D10").Formula = "=B10*C10"
$ws2.Range("E10").Value = 2
$ws2.Range("F10").Formula = "=E10*D10"

# row 11 - o1
$ws2.Range("B11").Formula = "=conditions!F6/60"
$ws2.Range("C11").Value = 2
$ws2.Range("D11").Formula = "=B11*C11"
$ws2.Range("E11").Value = 2
$ws2.Range("F11").Formula = "=E11*D11"

# row 12 - o2
$ws2.Range("B12").Formula = "=conditions!F7/60"
$ws2.Range("C12").Value = 2
$ws2.Range("D12").Formula = "=B12*C12"
$ws2.Range("E12").Value = 2
$ws2.Range("F12").Formula = "=E12*D12"

# summary row 13 (total over 2 days per response type)
$ws2.Range("I13").Formula = "=+F7"
$ws2.Range("J13").Formula = "=F8"
$ws2.Range("K13").Formula = "=SUM(F9:F10)"
$ws2.Range("L13").Formula = "=SUM(F11:F12)"

# total/day row 14
$ws2.Range("D14").Formula = "=SUM(D7:D12)"

# column H is a bit wider so the "total 2 days" label fits
$ws2.Columns.Item(8).ColumnWidth = 10.17

# leave the cursor on the README sheet where the new table was built, then
# return focus to the conditions sheet (which stays the active tab)
$ws2.Select()
$ws2.Range("C15").Select() | Out-Null
$ws1.Select()
